$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextCell($cell, $val) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextCell 'D2' '313.06'
Set-TextCell 'E2' '0.90%'
Set-TextCell 'G2' '16'
Set-TextCell 'D3' '37.85'
Set-TextCell 'E3' '-0.28%'
Set-TextCell 'G3' '16'
Set-TextCell 'E4' '0.88%'
Set-TextCell 'G4' '16'
Set-TextCell 'D5' '0.07918'
Set-TextCell 'E5' '0.67%'
Set-TextCell 'G5' '16'
Set-TextCell 'B6' 'FTXToken'
Set-TextCell 'C6' 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextCell 'D6' '1.917'
Set-TextCell 'E6' '-2.56%'
Set-TextCell 'G6' '16'
Set-TextCell 'B7' 'KuCoinToken'
Set-TextCell 'C7' 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextCell 'D7' '8.257'
Set-TextCell 'E7' '-0.39%'
Set-TextCell 'G7' '16'
Set-TextCell 'B8' 'MXToken'
Set-TextCell 'C8' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell 'D8' '0.9272'
Set-TextCell 'E8' '-0.45%'
Set-TextCell 'G8' '16'
Set-TextCell 'B9' 'LiechtensteinCryptoassetsExchange'
Set-TextCell 'C9' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextCell 'D9' '0.1224'
Set-TextCell 'E9' '-10.59%'
Set-TextCell 'G9' '16'
Set-TextCell 'B10' 'WazirX'
Set-TextCell 'C10' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextCell 'D10' '0.1921'
Set-TextCell 'E10' '-3.67%'
Set-TextCell 'G10' '16'
Set-TextCell 'B11' 'MandalaExchangeToken'
Set-TextCell 'C11' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextCell 'D11' '0.09117'
Set-TextCell 'E11' '2.19%'
Set-TextCell 'G11' '16'
Set-TextCell 'B12' 'BitrueCoin'
Set-TextCell 'C12' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextCell 'D12' '0.03334'
Set-TextCell 'E12' '-2.89%'
Set-TextCell 'G12' '16'
Set-TextCell 'B13' 'BitMartToken'
Set-TextCell 'C13' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextCell 'D13' '0.09630'
Set-TextCell 'E13' '-0.98%'
Set-TextCell 'G13' '16'
Set-TextCell 'B14' 'BitForexToken'
Set-TextCell 'C14' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextCell 'D14' '0.001375'
Set-TextCell 'E14' '-0.94%'
Set-TextCell 'G14' '16'
Set-TextCell 'B15' 'TigerCash'
Set-TextCell 'C15' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextCell 'D15' '0.005741'
Set-TextCell 'E15' '-6.56%'
Set-TextCell 'G15' '16'
Set-TextCell 'B16' 'LEO'
Set-TextCell 'C16' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextCell 'D16' '3.532'
Set-TextCell 'E16' '-1.14%'
Set-TextCell 'G16' '16'
Set-TextCell 'B17' 'GateToken'
Set-TextCell 'C17' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextCell 'D17' '4.411'
Set-TextCell 'E17' '0.90%'
Set-TextCell 'G17' '16'
Set-TextCell 'D18' '3.105'
Set-TextCell 'E18' '1.18%'
Set-TextCell 'G18' '16'
Set-TextCell 'D19' '0.3451'
Set-TextCell 'E19' '-0.51%'
Set-TextCell 'G19' '16'
Set-TextCell 'D20' '5.288'
Set-TextCell 'E20' '5.77%'
Set-TextCell 'G20' '16'
Set-TextCell 'D21' '0.1272'
Set-TextCell 'E21' '-1.91%'
Set-TextCell 'G21' '16'
Set-TextCell 'D22' '0.2590'
Set-TextCell 'E22' '4.12%'
Set-TextCell 'G22' '16'
Set-TextCell 'G23' '16'
Set-TextCell 'D24' '0.04380'
Set-TextCell 'E24' '1.35%'
Set-TextCell 'G24' '16'
Set-TextCell 'D25' '0.001253'
Set-TextCell 'E25' '2.92%'
Set-TextCell 'G25' '16'
Set-TextCell 'D26' '0.004297'
Set-TextCell 'E26' '-5.85%'
Set-TextCell 'G26' '16'
Set-TextCell 'E27' '-9.63%'
Set-TextCell 'G27' '16'
Set-TextCell 'G28' '16'
Set-TextCell 'G29' '16'
Set-TextCell 'G30' '16'
Set-TextCell 'G31' '16'
Set-TextCell 'G32' '16'
Set-TextCell 'G33' '16'
Set-TextCell 'G34' '16'
Set-TextCell 'G35' '16'
Set-TextCell 'G36' '16'
Set-TextCell 'G37' '16'
Set-TextCell 'G38' '16'
Set-TextCell 'D39' '0.02113'
Set-TextCell 'E39' '-6.98%'
Set-TextCell 'G39' '16'
Set-TextCell 'D40' '0.05115'
Set-TextCell 'E40' '1.17%'
Set-TextCell 'G40' '16'
Set-TextCell 'D41' '0.007643'
Set-TextCell 'E41' '2.16%'
Set-TextCell 'G41' '16'
Set-TextCell 'D42' '0.009133'
Set-TextCell 'E42' '-7.99%'
Set-TextCell 'G42' '16'
Set-TextCell 'D43' '0.1359'
Set-TextCell 'E43' '0.24%'
Set-TextCell 'G43' '16'
Set-TextCell 'D44' '0.001989'
Set-TextCell 'E44' '0.45%'
Set-TextCell 'G44' '16'
Set-TextCell 'D45' '0.008623'
Set-TextCell 'E45' '-1.92%'
Set-TextCell 'G45' '16'
Set-TextCell 'D46' '0.00006698'
Set-TextCell 'E46' '1.77%'
Set-TextCell 'G46' '16'
Set-TextCell 'E47' '0.03%'
Set-TextCell 'G47' '16'
Set-TextCell 'D48' '0.002913'
Set-TextCell 'E48' '-2.88%'
Set-TextCell 'G48' '16'
Set-TextCell 'E49' '-0.03%'
Set-TextCell 'G49' '16'
Set-TextCell 'D50' '0.00002100'
Set-TextCell 'E50' '0.03%'
Set-TextCell 'G50' '16'
Set-TextCell 'D51' '0.0002000'
Set-TextCell 'E51' '0.03%'
Set-TextCell 'G51' '16'
